# Automatische test-sync: 2025-08-04 20:09:50
$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 7 with the new test-mail entry ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(7, 1).Value = "Kun jij dit even regelen?"
$logs.Cells.Item(7, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(7, 3).Value = "Testmail #1: Kun jij dit even regelen?"
$logs.Cells.Item(7, 4).Value = "Planning / Afspraak"
$logs.Cells.Item(7, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Cells.Item(7, 6).Value = "2025-08-04 20:09:20"
$logs.Cells.Item(7, 7).Value = "Ja"
$logs.Cells.Item(7, 8).Value = "Ja"
$logs.Cells.Item(7, 9).Value = "Nee"
$logs.Cells.Item(7, 10).Value = "Nee"

# Extend the conditional formatting ranges (D/G/H/I/J) from row 6 to row 7
# while keeping the existing rules (priorities / dxfIds) intact.
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $logs.Range($col + "2:" + $col + "6")
    $newRange = $logs.Range($col + "2:" + $col + "7")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard sheet: append summary row for the new "Planning / Afspraak" category ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(4, 1).Value = "Planning / Afspraak"
$dash.Cells.Item(4, 2).Value = 1

# --- Chart: extend category/value series ranges from row 3 to row 4 ---
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$4,Dashboard!`$B`$2:`$B`$4,1)"
